$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 98
$ws.Cells.Item(98, 8).Value = 1797.75  # was 1848.7407
$ws.Cells.Item(98, 9).Value = 1664.875  # was 1718.9565
$ws.Cells.Item(98, 11).Value = 1664.875  # was 1718.9565
$ws.Cells.Item(98, 13).Value = -166.875  # was -220.9565
# Row 103
$ws.Cells.Item(103, 8).Value = 20834290  # was 21740098
$ws.Cells.Item(103, 9).Value = 751.3125  # was 755.4666999999999
$ws.Cells.Item(103, 11).Value = 2253.9375  # was 2266.4001
$ws.Cells.Item(103, 13).Value = -1667.9375  # was -1680.4001
# Row 122
$ws.Cells.Item(122, 8).Value = 1797.75  # was 1848.7407
$ws.Cells.Item(122, 9).Value = 1664.875  # was 1718.9565
$ws.Cells.Item(122, 11).Value = 4994.625  # was 5156.8695
$ws.Cells.Item(122, 13).Value = -2544.625  # was -2706.8695
# Row 129
$ws.Cells.Item(129, 8).Value = 3353.2307  # was 3821.7727
$ws.Cells.Item(129, 9).Value = 797.6  # was 859.0833
$ws.Cells.Item(129, 10).Value = 6838.1816  # was 7377
$ws.Cells.Item(129, 11).Value = 2392.8  # was 2577.2499
$ws.Cells.Item(129, 12).Value = 20514.5448  # was 22131
$ws.Cells.Item(129, 13).Value = 2607.2  # was 2422.7501
$ws.Cells.Item(129, 14).Value = -30514.5448  # was -32131
# Row 132
$ws.Cells.Item(132, 8).Value = 6980.3105  # was 7349.778
$ws.Cells.Item(132, 9).Value = 3727.25  # was 3975.0715
$ws.Cells.Item(132, 11).Value = 11181.75  # was 11925.2145
$ws.Cells.Item(132, 13).Value = -8651.75  # was -9395.2145
# Row 137
$ws.Cells.Item(137, 8).Value = 2092.2173  # was 1891.6
$ws.Cells.Item(137, 9).Value = 2079.0527  # was 1877.24
$ws.Cells.Item(137, 10).Value = 2154.75  # was 1963.4
$ws.Cells.Item(137, 11).Value = 6237.158100000001  # was 5631.72
$ws.Cells.Item(137, 12).Value = 6464.25  # was 5890.200000000001
$ws.Cells.Item(137, 13).Value = -3687.158100000001  # was -3081.72
$ws.Cells.Item(137, 14).Value = -11564.25  # was -10990.2
# Row 141
$ws.Cells.Item(141, 8).Value = 35755384  # was 38505260
$ws.Cells.Item(141, 9).Value = 71435880  # was 83340690
$ws.Cells.Item(141, 11).Value = 214307640  # was 250022070
$ws.Cells.Item(141, 13).Value = -214302460  # was -250016890

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 3580.5173  # was 3707.0356
$ws.Cells.Item(32, 9).Value = 2419.8333  # was 2511.4424
$ws.Cells.Item(32, 11).Value = 2419.8333  # was 2511.4424
$ws.Cells.Item(32, 13).Value = -2132.8333  # was -2224.4424
# Row 45
$ws.Cells.Item(45, 8).Value = 4470.6665  # was 2276.8333
$ws.Cells.Item(45, 9).Value = 2706  # was 1132.2
$ws.Cells.Item(45, 11).Value = 2706  # was 1132.2
$ws.Cells.Item(45, 13).Value = -2329  # was -755.2
# Row 61
$ws.Cells.Item(61, 8).Value = 32858028  # was 11580509
$ws.Cells.Item(61, 9).Value = 42001240  # was 13335151
$ws.Cells.Item(61, 10).Value = 10000000  # was 5000600
$ws.Cells.Item(61, 11).Value = 42001240  # was 13335151
$ws.Cells.Item(61, 12).Value = 10000000  # was 5000600
$ws.Cells.Item(61, 13).Value = -42001028  # was -13334939
$ws.Cells.Item(61, 14).Value = -10000424  # was -5001024
# Row 74
$ws.Cells.Item(74, 8).Value = 2436.6  # was 2609.4546
$ws.Cells.Item(74, 9).Value = 1455.85  # was 1506.4706
$ws.Cells.Item(74, 11).Value = 1455.85  # was 1506.4706
$ws.Cells.Item(74, 13).Value = -581.8499999999999  # was -632.4706000000001
# Row 77
$ws.Cells.Item(77, 8).Value = 2436.6  # was 2609.4546
$ws.Cells.Item(77, 9).Value = 1455.85  # was 1506.4706
$ws.Cells.Item(77, 11).Value = 7279.25  # was 7532.353000000001
$ws.Cells.Item(77, 13).Value = -2911.25  # was -3164.353000000001
# Row 97
$ws.Cells.Item(97, 8).Value = 1877.1538  # was 1910.2727
$ws.Cells.Item(97, 9).Value = 1854.8182  # was 1910.2727
$ws.Cells.Item(97, 10).Value = 2000  # was 0
$ws.Cells.Item(97, 11).Value = 1854.8182  # was 1910.2727
$ws.Cells.Item(97, 12).Value = 2000  # was 0
$ws.Cells.Item(97, 13).Value = -1358.8182  # was -1414.2727
$ws.Cells.Item(97, 14).Value = -2992  # was None
# Row 132
$ws.Cells.Item(132, 8).Value = 8340672  # was 3851316.8
$ws.Cells.Item(132, 9).Value = 7852.2856  # was 4929.778
$ws.Cells.Item(132, 10).Value = 20006620  # was 12505688
$ws.Cells.Item(132, 11).Value = 23556.8568  # was 14789.334
$ws.Cells.Item(132, 12).Value = 60019860  # was 37517064
$ws.Cells.Item(132, 13).Value = -21026.8568  # was -12259.334
$ws.Cells.Item(132, 14).Value = -60024920  # was -37522124
# Row 136
$ws.Cells.Item(136, 8).Value = 32858028  # was 11580509
$ws.Cells.Item(136, 9).Value = 42001240  # was 13335151
$ws.Cells.Item(136, 10).Value = 10000000  # was 5000600
$ws.Cells.Item(136, 11).Value = 126003720  # was 40005453
$ws.Cells.Item(136, 12).Value = 30000000  # was 15001800
$ws.Cells.Item(136, 13).Value = -126001170  # was -40002903
$ws.Cells.Item(136, 14).Value = -30005100  # was -15006900

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Cells.Item(134, 8).Value = 10003068  # was 4002035.2
$ws.Cells.Item(134, 9).Value = 3241.4285  # was 2118.476
$ws.Cells.Item(134, 10).Value = 33335996  # was 25001598
$ws.Cells.Item(134, 11).Value = 9724.2855  # was 6355.428
$ws.Cells.Item(134, 12).Value = 100007988  # was 75004794
$ws.Cells.Item(134, 13).Value = -7189.2855  # was -3820.428
$ws.Cells.Item(134, 14).Value = -100013058  # was -75009864

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 12
$ws.Cells.Item(12, 8).Value = 4495  # was 4496
$ws.Cells.Item(12, 9).Value = 4495  # was 4496
$ws.Cells.Item(12, 11).Value = 4495  # was 4496
$ws.Cells.Item(12, 13).Value = -4325  # was -4326
# Row 16
$ws.Cells.Item(16, 8).Value = 11113767  # was 16669167
$ws.Cells.Item(16, 9).Value = 12502612  # was 20002400
$ws.Cells.Item(16, 11).Value = 12502612  # was 20002400
$ws.Cells.Item(16, 13).Value = -12502325  # was -20002113
# Row 31
$ws.Cells.Item(31, 8).Value = 27780870  # was 27030128
$ws.Cells.Item(31, 9).Value = 43481132  # was 41669560
$ws.Cells.Item(31, 11).Value = 43481132  # was 41669560
$ws.Cells.Item(31, 13).Value = -43480837  # was -41669265
# Row 32
$ws.Cells.Item(32, 8).Value = 2730.5  # was 3327.75
$ws.Cells.Item(32, 9).Value = 2730.5  # was 3327.75
$ws.Cells.Item(32, 11).Value = 2730.5  # was 3327.75
$ws.Cells.Item(32, 13).Value = -2414.5  # was -3011.75
# Row 34
$ws.Cells.Item(34, 8).Value = 27780870  # was 27030128
$ws.Cells.Item(34, 9).Value = 43481132  # was 41669560
$ws.Cells.Item(34, 11).Value = 43481132  # was 41669560
$ws.Cells.Item(34, 13).Value = -43480930  # was -41669358
# Row 99
$ws.Cells.Item(99, 8).Value = 42656  # was 66999.60000000001
$ws.Cells.Item(99, 9).Value = 13167.333  # was 14999
$ws.Cells.Item(99, 10).Value = 101633.336  # was 101666.664
$ws.Cells.Item(99, 11).Value = 13167.333  # was 14999
$ws.Cells.Item(99, 12).Value = 101633.336  # was 101666.664
$ws.Cells.Item(99, 13).Value = -11669.333  # was -13501
$ws.Cells.Item(99, 14).Value = -104629.336  # was -104662.664
# Row 113
$ws.Cells.Item(113, 8).Value = 11113767  # was 16669167
$ws.Cells.Item(113, 9).Value = 12502612  # was 20002400
$ws.Cells.Item(113, 11).Value = 12502612  # was 20002400
$ws.Cells.Item(113, 13).Value = -12500442  # was -20000230
# Row 118
$ws.Cells.Item(118, 8).Value = 99999  # was 0
$ws.Cells.Item(118, 10).Value = 99999  # was 0
$ws.Cells.Item(118, 12).Value = 99999  # was 0
$ws.Cells.Item(118, 14).Value = -103313  # was None
# Row 122
$ws.Cells.Item(122, 8).Value = 6240  # was 4110
$ws.Cells.Item(122, 9).Value = 10000  # was 3584
$ws.Cells.Item(122, 11).Value = 30000  # was 10752
$ws.Cells.Item(122, 13).Value = -27550  # was -8302
# Row 126
$ws.Cells.Item(126, 8).Value = 42656  # was 66999.60000000001
$ws.Cells.Item(126, 9).Value = 13167.333  # was 14999
$ws.Cells.Item(126, 10).Value = 101633.336  # was 101666.664
$ws.Cells.Item(126, 11).Value = 39501.999  # was 44997
$ws.Cells.Item(126, 12).Value = 304900.008  # was 304999.992
$ws.Cells.Item(126, 13).Value = -37031.999  # was -42527
$ws.Cells.Item(126, 14).Value = -309840.008  # was -309939.992
# Row 132
$ws.Cells.Item(132, 8).Value = 2465.923  # was 2544.5405
$ws.Cells.Item(132, 9).Value = 2398.7878  # was 2442.125
$ws.Cells.Item(132, 10).Value = 2835.1667  # was 3200
$ws.Cells.Item(132, 11).Value = 7196.3634  # was 7326.375
$ws.Cells.Item(132, 12).Value = 8505.500100000001  # was 9600
$ws.Cells.Item(132, 13).Value = -4666.3634  # was -4796.375
$ws.Cells.Item(132, 14).Value = -13565.5001  # was -14660
# Row 134
$ws.Cells.Item(134, 8).Value = 0  # was 2859.5264
$ws.Cells.Item(134, 9).Value = 0  # was 2645.875
$ws.Cells.Item(134, 10).Value = 0  # was 3999
$ws.Cells.Item(134, 11).Value = 0  # was 7937.625
$ws.Cells.Item(134, 12).ClearContents()  # was 11997
$ws.Cells.Item(134, 13).ClearContents()  # was -5402.625
$ws.Cells.Item(134, 14).Value = 0  # was -17067
# Row 141
$ws.Cells.Item(141, 8).Value = 261666.67  # was 268333.34
$ws.Cells.Item(141, 10).Value = 261666.67  # was 268333.34
$ws.Cells.Item(141, 12).Value = 261666.67  # was 268333.34
$ws.Cells.Item(141, 14).Value = -272026.67  # was -278693.34

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Cells.Item(7, 8).Value = 571  # was 656.7143
$ws.Cells.Item(7, 9).Value = 639.6  # was 649.5
$ws.Cells.Item(7, 10).Value = 399.5  # was 700
$ws.Cells.Item(7, 11).Value = 1918.8  # was 1948.5
$ws.Cells.Item(7, 12).Value = 1198.5  # was 2100
$ws.Cells.Item(7, 13).Value = -1806.8  # was -1836.5
$ws.Cells.Item(7, 14).Value = -1422.5  # was -2324
# Row 125
$ws.Cells.Item(125, 8).Value = 9999  # was 21032.727
$ws.Cells.Item(125, 9).Value = 0  # was 6514.5
$ws.Cells.Item(125, 10).Value = 9999  # was 24259
$ws.Cells.Item(125, 11).Value = 0  # was 19543.5
$ws.Cells.Item(125, 12).ClearContents()  # was 72777
$ws.Cells.Item(125, 13).Value = 29997  # was -14623.5
$ws.Cells.Item(125, 14).Value = -39837  # was -82617
# Row 131
$ws.Cells.Item(131, 8).Value = 5151.9443  # was 5237.4116
$ws.Cells.Item(131, 9).Value = 3351.5  # was 3247.5386
$ws.Cells.Item(131, 10).Value = 8752.833000000001  # was 11704.5
$ws.Cells.Item(131, 11).Value = 10054.5  # was 9742.6158
$ws.Cells.Item(131, 12).Value = 26258.499  # was 35113.5
$ws.Cells.Item(131, 13).Value = -5014.5  # was -4702.6158
$ws.Cells.Item(131, 14).Value = -36338.499  # was -45193.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 28
$ws.Cells.Item(28, 8).Value = 2000  # was 0
$ws.Cells.Item(28, 9).Value = 2000  # was 0
$ws.Cells.Item(28, 11).Value = 2000  # was 0
$ws.Cells.Item(28, 13).Value = -1808  # was None
# Row 70
$ws.Cells.Item(70, 8).Value = 5710.543  # was 5814.3335
$ws.Cells.Item(70, 9).Value = 6484.857  # was 7479.6
$ws.Cells.Item(70, 11).Value = 6484.857  # was 7479.6
$ws.Cells.Item(70, 13).Value = -6214.857  # was -7209.6
# Row 73
$ws.Cells.Item(73, 8).Value = 5710.543  # was 5814.3335
$ws.Cells.Item(73, 9).Value = 6484.857  # was 7479.6
$ws.Cells.Item(73, 11).Value = 6484.857  # was 7479.6
$ws.Cells.Item(73, 13).Value = -5548.857  # was -6543.6
# Row 80
$ws.Cells.Item(80, 8).Value = 2055.3076  # was 2204.625
$ws.Cells.Item(80, 9).Value = 1778  # was 1872.3
$ws.Cells.Item(80, 10).Value = 2499  # was 2758.5
$ws.Cells.Item(80, 11).Value = 1778  # was 1872.3
$ws.Cells.Item(80, 12).Value = 2499  # was 2758.5
$ws.Cells.Item(80, 13).Value = -780  # was -874.3
$ws.Cells.Item(80, 14).Value = -4495  # was -4754.5
# Row 83
$ws.Cells.Item(83, 8).Value = 2055.3076  # was 2204.625
$ws.Cells.Item(83, 9).Value = 1778  # was 1872.3
$ws.Cells.Item(83, 10).Value = 2499  # was 2758.5
$ws.Cells.Item(83, 11).Value = 8890  # was 9361.5
$ws.Cells.Item(83, 12).Value = 12495  # was 13792.5
$ws.Cells.Item(83, 13).Value = -3898  # was -4369.5
$ws.Cells.Item(83, 14).Value = -22479  # was -23776.5
# Row 92
$ws.Cells.Item(92, 8).Value = 51748.668  # was 55082
$ws.Cells.Item(92, 10).Value = 51748.668  # was 55082
$ws.Cells.Item(92, 12).Value = 51748.668  # was 55082
$ws.Cells.Item(92, 14).Value = -55492.668  # was -58826
# Row 102
$ws.Cells.Item(102, 8).Value = 3210.7727  # was 3402.6316
$ws.Cells.Item(102, 9).Value = 3263.7778  # was 3517.4
$ws.Cells.Item(102, 11).Value = 3263.7778  # was 3517.4
$ws.Cells.Item(102, 13).Value = -1641.7778  # was -1895.4
# Row 126
$ws.Cells.Item(126, 8).Value = 6683  # was 2838.6
$ws.Cells.Item(126, 9).Value = 0  # was 1599.5
$ws.Cells.Item(126, 10).Value = 6683  # was 3664.6667
$ws.Cells.Item(126, 11).Value = 0  # was 4798.5
$ws.Cells.Item(126, 12).ClearContents()  # was 10994.0001
$ws.Cells.Item(126, 13).Value = 20049  # was -2328.5
$ws.Cells.Item(126, 14).Value = -24989  # was -15934.0001
# Row 132
$ws.Cells.Item(132, 8).Value = 100000000  # was 9094909
$ws.Cells.Item(132, 9).Value = 0  # was 4333.5557
$ws.Cells.Item(132, 10).Value = 100000000  # was 50002500
$ws.Cells.Item(132, 11).Value = 0  # was 13000.6671
$ws.Cells.Item(132, 12).ClearContents()  # was 150007500
$ws.Cells.Item(132, 13).Value = 300000000  # was -10470.6671
$ws.Cells.Item(132, 14).Value = -300005060  # was -150012560
# Row 136
$ws.Cells.Item(136, 8).Value = 64882.832  # was 62027.855
$ws.Cells.Item(136, 10).Value = 64882.832  # was 62027.855
$ws.Cells.Item(136, 12).Value = 194648.496  # was 186083.565
$ws.Cells.Item(136, 14).Value = -199748.496  # was -191183.565

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Cells.Item(40, 8).Value = 3353.3809  # was 3353.4285
$ws.Cells.Item(40, 9).Value = 3441.2778  # was 3441.3333
$ws.Cells.Item(40, 11).Value = 3441.2778  # was 3441.3333
$ws.Cells.Item(40, 13).Value = -3305.2778  # was -3305.3333
# Row 93
$ws.Cells.Item(93, 8).Value = 1987599.8  # was 1855152.8
$ws.Cells.Item(93, 9).Value = 1984.7693  # was 1839.4667
$ws.Cells.Item(93, 11).Value = 1984.7693  # was 1839.4667
$ws.Cells.Item(93, 13).Value = -736.7692999999999  # was -591.4666999999999
# Row 132
$ws.Cells.Item(132, 8).Value = 4239.263  # was 4026.1428
$ws.Cells.Item(132, 9).Value = 3042.2307  # was 2903.4666
$ws.Cells.Item(132, 11).Value = 9126.6921  # was 8710.399800000001
$ws.Cells.Item(132, 13).Value = -6596.6921  # was -6180.399800000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 41
$ws.Cells.Item(41, 8).Value = 18678  # was 23247.25
$ws.Cells.Item(41, 10).Value = 16598.75  # was 21998
$ws.Cells.Item(41, 12).Value = 16598.75  # was 21998
$ws.Cells.Item(41, 14).Value = -17378.75  # was -22778
# Row 95
$ws.Cells.Item(95, 8).Value = 99999  # was 99498.5
$ws.Cells.Item(95, 10).Value = 99999  # was 99498.5
$ws.Cells.Item(95, 12).Value = 99999  # was 99498.5
$ws.Cells.Item(95, 14).Value = -105491  # was -104990.5
# Row 109
$ws.Cells.Item(109, 8).Value = 99999  # was 57499.5
$ws.Cells.Item(109, 9).Value = 0  # was 15000
$ws.Cells.Item(109, 11).Value = 0  # was 15000
$ws.Cells.Item(109, 13).ClearContents()  # was -13613
# Row 126
$ws.Cells.Item(126, 8).Value = 3698.7144  # was 3840.15
$ws.Cells.Item(126, 9).Value = 4115.222  # was 4306.1177
$ws.Cells.Item(126, 11).Value = 12345.666  # was 12918.3531
$ws.Cells.Item(126, 13).Value = -9875.665999999999  # was -10448.3531
# Row 132
$ws.Cells.Item(132, 8).Value = 5002250  # was 715491.0600000001
$ws.Cells.Item(132, 9).Value = 4500  # was 1302.0834
$ws.Cells.Item(132, 10).Value = 10000000  # was 5000625
$ws.Cells.Item(132, 11).Value = 13500  # was 3906.2502
$ws.Cells.Item(132, 12).Value = 30000000  # was 15001875
$ws.Cells.Item(132, 13).Value = -10970  # was -1376.2502
$ws.Cells.Item(132, 14).Value = -30005060  # was -15006935
